# Daily attendance processing - 2026-01-01 07:58:19
# For every row in the "Recorded By" column (G), if the comma-separated
# list of recorders contains the exact token "System", reverse the order
# of the whole comma-separated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        if (($parts.Length -gt 1) -and ($parts -contains "System")) {
            $n = $parts.Length
            $reversed = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $newVal = $reversed -join ", "
            $cell.Value2 = $newVal
        }
    }
}
